# Update to CrowdFlo import process
# Re-label the "Question N" header columns (M1:AB1) to match the new
# CrowdFlo question ordering, and bump the header/data row heights from
# 18.75 to 19.5 (18.75pt -> 19.5pt, i.e. 25px -> 26px row height).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-point each "Question N" header cell at its new question number.
# (Row 2 answer values are untouched - only the header labels move.)
$ws.Range("M1").Value  = "Question 16"
$ws.Range("N1").Value  = "Question 6"
$ws.Range("O1").Value  = "Question 8"
$ws.Range("P1").Value  = "Question 15"
$ws.Range("Q1").Value  = "Question 14"
$ws.Range("R1").Value  = "Question 3"
$ws.Range("S1").Value  = "Question 9"
$ws.Range("T1").Value  = "Question 13"
$ws.Range("U1").Value  = "Question 2"
$ws.Range("V1").Value  = "Question 7"
$ws.Range("W1").Value  = "Question 4"
$ws.Range("X1").Value  = "Question 12"
$ws.Range("Y1").Value  = "Question 10"
$ws.Range("Z1").Value  = "Question 8"
$ws.Range("AA1").Value = "Question 11"
$ws.Range("AB1").Value = "Question 1"

# Row heights: both data rows grow slightly (18.75 -> 19.5).
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
